$wb = $excel.ActiveWorkbook

# Rename sheets (new participant-generation timestamps)
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911336104875"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911368724408"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911368734398"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911369302914"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911370191853"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911335680077.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911335858073.csv"
$ws1.Range("B4").Value = "go_stims-1650291133586446.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911336094878.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650291135881013.csv"
$ws2.Range("B3").Value = "TB-16502911368497703.csv"
$ws2.Range("B4").Value = "OB-16502911354647012.csv"
$ws2.Range("B5").Value = "TB-16502911367419972.csv"
$ws2.Range("B6").Value = "ZB-match_6-1650291133768858.csv"
$ws2.Range("B7").Value = "TB-16502911366693418.csv"
$ws2.Range("B8").Value = "OB-16502911344423692.csv"
$ws2.Range("B9").Value = "ZB-match_7-16502911336910946.csv"
$ws2.Range("B10").Value = "ZB-match_6-16502911339603798.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911368883786.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911368754456.csv"
$ws4.Range("B4").Value = "MM_stims-16502911369141219.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911368893423.csv"
$ws4.Range("B6").Value = "MM_stims-16502911369293287.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911369151266.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502911369769409.csv"
$ws5.Range("B3").Value = "SAT_stims-16502911369352934.csv"
$ws5.Range("B4").Value = "SAT_stims-1650291136960761.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502911370031323.csv"
